$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "mre_costs" worksheet right after "pv_costs" (and before
#    "battery_costs"), to support multiple types of renewable energy (RE).
# ---------------------------------------------------------------------------
$pv = $wb.Worksheets.Item("pv_costs")
$mre = $wb.Worksheets.Add($null, $pv)
$mre.Name = "mre_costs"
$mre.Range("A1").Value = "Tidal"
$mre.Range("A2").Value = 10000
$mre.Activate()
$mre.Range("A2").Select()

# ---------------------------------------------------------------------------
# 2. pv_costs: rows 2-5 in column A lose their (no-op) alignment style.
# ---------------------------------------------------------------------------
$pv.Range("A2:A5").ClearFormats()

# ---------------------------------------------------------------------------
# 3. battery_costs: C3 picks up the same centered style already used by the
#    rest of the table (style de-duplication).
# ---------------------------------------------------------------------------
$battery = $wb.Worksheets.Item("battery_costs")
$battery.Range("C3").HorizontalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------------------
# 4. om_costs: add the new "Tidal" / "$/turbine-yr" / 100 column, clean up
#    the redundant font-only styles on the existing cells, and make this the
#    active sheet/tab (as in the saved workbook view state).
# ---------------------------------------------------------------------------
$om = $wb.Worksheets.Item("om_costs")

# New column I: Tidal unit cost row.
$om.Range("I1").Value = "Tidal"
$om.Range("I2").Value = "$/turbine-yr"
$om.Range("I3").Value = 100
$om.Columns.Item(9).ColumnWidth = 14.166666666666666

$om.Range("I2").HorizontalAlignment = -4108        # xlCenter
$om.Range("I2").VerticalAlignment = -4108          # xlCenter
$om.Range("I3").HorizontalAlignment = -4108        # xlCenter

# Drop the redundant "applyFont" flag baked into the header/unit rows -
# clear then re-apply just the alignment that should remain.
$om.Range("B1:D1").ClearFormats()
$om.Range("B1:D1").HorizontalAlignment = -4108      # xlCenter

$om.Range("E1:H1").ClearFormats()

$om.Range("D2").ClearFormats()
$om.Range("D2").HorizontalAlignment = -4108         # xlCenter

$om.Range("B2:H2").ClearFormats()
$om.Range("B2:H2").HorizontalAlignment = -4108      # xlCenter
$om.Range("B2:H2").VerticalAlignment = -4108        # xlCenter

$om.Range("B3:H3").ClearFormats()
$om.Range("B3:H3").HorizontalAlignment = -4108      # xlCenter

$om.Activate()
$om.Range("J3").Select()
